$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResults")

$ws.Range("A2").Value = "iAU_TC_ID_107"
$ws.Range("B2").Value = "@RegressionA Validation of Blueprints list page"
$ws.Range("C2").Value = "passed"
